$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.223.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.427.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.57%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.516'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.47%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  +7.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.806.93'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.439.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.126.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("E20").Value = '  -0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0919'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  -12.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.43%  '
$ws.Range("E33").Value = '  +4.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.20'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.02%  '
$ws.Range("E35").Value = '  +0.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0766'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.90%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  +2.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.938.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("E47").Value = '  +3.11%  '
$ws.Range("E48").Value = '  +16.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.55%  '
$ws.Range("E51").Value = '  +2.12%  '
